# Adyen "credit fees" test file migration (12.0):
# The Gross/Net Currency columns (K and O) on the "Data" sheet were
# re-sampled from EUR to USD for every transaction row except the one
# GBP outlier in row 28 (column K only - its Net Currency in column O
# still becomes USD like everything else).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Gross Currency (column K): rows 5-27 and 29-42 -> USD.
# Row 28 keeps its original "GBP" value and is intentionally skipped.
$ws.Range("K5:K27").Value = "USD"
$ws.Range("K29:K42").Value = "USD"

# Net Currency (column O): rows 5-44 (including row 28) -> USD.
$ws.Range("O5:O44").Value = "USD"

# Restore the active selection to what it was left at when the file was
# last saved (P38 on the Data sheet).
$ws.Range("P38").Select() | Out-Null
